$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell -> new value for this "Updated symbol list" data refresh.
# Price-like columns (numeric-looking text) need NumberFormat "@" forced
# first so Excel stores them as text (matching the source data, which
# keeps trailing zeros / fixed precision) instead of coercing to Number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "258.02"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.79"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.160"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06070"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.732"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.449"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7974"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1570"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08059"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03361"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03089"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09304"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.896"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001705"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04841"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006158"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006202"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001100"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003379"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.690"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.263"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003020"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04574"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007110"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003906"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1112"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009998"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.002974"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005922"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7510"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06784"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001502"
